$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '98.472.21'
$ws.Range('E2').Value = '  +4.44%  '
$ws.Range('D3').Value = '3.362.26'
$ws.Range('E3').Value = '  +9.53%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''257.04'
$ws.Range('E5').Value = '  +9.09%  '
$ws.Range('D6').Value = '''622.61'
$ws.Range('E6').Value = '  +2.73%  '
$ws.Range('D7').Value = '''1.24'
$ws.Range('E7').Value = '  +12.20%  '
$ws.Range('D8').Value = '''0.388'
$ws.Range('E8').Value = '  +2.51%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '3.362.70'
$ws.Range('E10').Value = '  +9.56%  '
$ws.Range('D11').Value = '''0.819'
$ws.Range('E11').Value = '  +2.42%  '
$ws.Range('E12').Value = '  +1.90%  '
$ws.Range('D13').Value = '98.142.81'
$ws.Range('E13').Value = '  +4.56%  '
$ws.Range('E14').Value = '  +6.56%  '
$ws.Range('E15').Value = '  +3.74%  '
$ws.Range('D16').Value = '3.997.08'
$ws.Range('E16').Value = '  +9.77%  '
$ws.Range('E17').Value = '  +4.20%  '
$ws.Range('D18').Value = '3.367.08'
$ws.Range('E18').Value = '  +9.75%  '
$ws.Range('D19').Value = '''3.63'
$ws.Range('E19').Value = '  +3.01%  '
$ws.Range('E20').Value = '  +4.74%  '
$ws.Range('D21').Value = '''486.35'
$ws.Range('E21').Value = '  +10.50%  '
$ws.Range('E22').Value = '  +3.76%  '
$ws.Range('D23').Value = '''0.0000207'
$ws.Range('D24').Value = '''9.27'
$ws.Range('E24').Value = '  +5.18%  '
$ws.Range('D25').Value = '''5.73'
$ws.Range('E25').Value = '  +4.53%  '
$ws.Range('D26').Value = '''88.33'
$ws.Range('E26').Value = '  +4.65%  '
$ws.Range('E27').Value = '  +2.59%  '
$ws.Range('E28').Value = '  +9.55%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '''0.261'
$ws.Range('E30').Value = '  +5.29%  '
$ws.Range('D31').Value = '''0.186'
$ws.Range('E31').Value = '  +4.63%  '
$ws.Range('D32').Value = '''0.127'
$ws.Range('E32').Value = '  +3.83%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').Value = '''9.23'
$ws.Range('E34').Value = '  +4.47%  '
$ws.Range('D35').Value = '''27.37'
$ws.Range('E35').Value = '  +7.72%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '''0.153'
$ws.Range('E36').Value = '  +0.99%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').Value = '''519.85'
$ws.Range('E37').Value = '  +7.97%  '
$ws.Range('D38').Value = '''7.37'
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('E39').Value = '  +4.60%  '
$ws.Range('D40').Value = '''24.81'
$ws.Range('E40').Value = '  +3.19%  '
$ws.Range('D41').Value = '''0.454'
$ws.Range('E41').Value = '  +5.11%  '
$ws.Range('B42').Value = 'MantraDAO'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D42').Value = '''3.75'
$ws.Range('E42').Value = '  -1.39%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').Value = '''1.28'
$ws.Range('E43').Value = '  +3.20%  '
$ws.Range('E44').Value = '  +7.28%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').Value = '''1.00'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').Value = '''0.781'
$ws.Range('E46').Value = '  +16.18%  '
$ws.Range('D47').Value = '''160.30'
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('E48').Value = '  +6.63%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = '''45.52'
$ws.Range('E49').Value = '  +4.58%  '
$ws.Range('E50').Value = '  +6.90%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').Value = '''4.53'
$ws.Range('E51').Value = '  +6.81%  '
